# Actualización automática de dinoe.xlsx desde Google Drive
#
# Updates the "productos" sheet: revises several "Precio Cliente Final en
# Soles" (column E) values and populates a new column G with updated
# prices for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("productos")
$ws.Activate()

# --- Column E (Precio Cliente Final en Soles) updates ---
$ws.Range("E28").Value = 33.2

$ws.Range("E85").Value = 19.16
$ws.Range("E86").Value = 34.35
$ws.Range("E87").Value = 7.72
$ws.Range("E88").Value = 13.14
$ws.Range("E90").Value = 137.39

$ws.Range("E97").Value = 3.78
$ws.Range("E98").Value = 3.78
$ws.Range("E99").Value = 3.86
$ws.Range("E100").Value = 3.86
$ws.Range("E101").Value = 3.86

# --- New column G values ---
$ws.Range("G121").Value = 29.5
$ws.Range("G122").Value = 33.5
$ws.Range("G123").Value = 32.5
$ws.Range("G124").Value = 28.5

$ws.Range("G157").Value = 30.0
$ws.Range("G158").Value = 33.5
$ws.Range("G159").Value = 32.5
$ws.Range("G160").Value = 28.5
